# Apply "Multiple Data Added in ORG Test" edit:
#  - Add 3 new rows (Tata/Energy, Windows/Communications, Mac/Electronics)
#    to the DataProviderOrg sheet, with Consolas font styling on column B.
#  - Widen column B on that sheet.
#  - Make DataProviderOrg the active/selected sheet (was Contact).

$wb = $excel.ActiveWorkbook

$contact = $wb.Worksheets.Item("Contact")
$ws = $wb.Worksheets.Item("DataProviderOrg")

# New data rows
$ws.Range("A7").Value = "Tata"
$ws.Range("B7").Value = "Energy"

$ws.Range("A8").Value = "Windows"
$ws.Range("B8").Value = "Communications"

$ws.Range("A9").Value = "Mac"
$ws.Range("B9").Value = "Electronics"

# B7 gets a non-bold Consolas font (derived from the existing bold Consolas
# style already used on the Contact sheet).
$contact.Range("F5").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Font.Bold = $false

# B8/B9 get a plain Consolas font with no explicit colour.
$ws.Range("B8").Font.Name = "Consolas"
$ws.Range("B8").Font.Size = 11
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# Widen column B to fit the new longer text.
$ws.Columns.Item(2).ColumnWidth = 17.5

# Make DataProviderOrg the active sheet/tab and select the last-entered cell.
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("B8").Select()
